$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G4").Value = "Miss Dina Nasr, Administrator"
$ws.Range("G5").Value = "Miss Dina Nasr, Administrator"
$ws.Range("G8").Value = "Miss Dina Nasr, Administrator"
$ws.Range("G17").Value = "Miss Dina Nasr, Administrator"
$ws.Range("G18").Value = "Miss Dina Nasr, Administrator"
$ws.Range("G19").Value = "Miss Dina Nasr, Administrator"
$ws.Range("G20").Value = "Miss Dina Nasr, Administrator"
$ws.Range("G21").Value = "Miss Dina Nasr, Administrator"
$ws.Range("G22").Value = "Miss Dina Nasr, Administrator"
$ws.Range("G23").Value = "Miss Dina Nasr, Administrator"
$ws.Range("G24").Value = "Miss Dina Nasr, Administrator"
$ws.Range("G26").Value = "Miss Dina Nasr, Administrator"
$ws.Range("G30").Value = "Miss Dina Nasr, Administrator"
$ws.Range("G31").Value = "Miss Dina Nasr, Administrator"
$ws.Range("G34").Value = "Miss Dina Nasr, Administrator"
$ws.Range("G43").Value = "Miss Dina Nasr, Administrator"
$ws.Range("G44").Value = "Miss Dina Nasr, Administrator"
$ws.Range("G45").Value = "Miss Dina Nasr, Administrator"
$ws.Range("G46").Value = "Miss Dina Nasr, Administrator"
$ws.Range("G47").Value = "Miss Dina Nasr, Administrator"
$ws.Range("G48").Value = "Miss Dina Nasr, Administrator"
$ws.Range("G49").Value = "Miss Dina Nasr, Administrator"
$ws.Range("G50").Value = "Miss Dina Nasr, Administrator"
$ws.Range("G52").Value = "Miss Dina Nasr, Administrator"
$ws.Range("G56").Value = "Miss Dina Nasr, Administrator"
$ws.Range("G57").Value = "Miss Dina Nasr, Administrator"
$ws.Range("G60").Value = "Miss Dina Nasr, Administrator"
$ws.Range("G69").Value = "Miss Dina Nasr, Administrator"
$ws.Range("G70").Value = "Miss Dina Nasr, Administrator"
$ws.Range("G71").Value = "Miss Dina Nasr, Administrator"
$ws.Range("G72").Value = "Miss Dina Nasr, Administrator"
$ws.Range("G73").Value = "Miss Dina Nasr, Administrator"
$ws.Range("G74").Value = "Miss Dina Nasr, Administrator"
$ws.Range("G75").Value = "Miss Dina Nasr, Administrator"
$ws.Range("G76").Value = "Miss Dina Nasr, Administrator"
$ws.Range("G78").Value = "Miss Dina Nasr, Administrator"
$ws.Range("G80").Value = "Miss Dina Nasr, Administrator"
$ws.Range("G81").Value = "Miss Dina Nasr, Administrator"
$ws.Range("G82").Value = "Miss Dina Nasr, Administrator"
$ws.Range("G93").Value = "Miss Dina Nasr, Administrator"
$ws.Range("G94").Value = "Miss Dina Nasr, Administrator"
$ws.Range("G96").Value = "Miss Dina Nasr, Administrator"
$ws.Range("G99").Value = "Miss Dina Nasr, Administrator"
$ws.Range("G101").Value = "Miss Dina Nasr, Administrator"
$ws.Range("G106").Value = "Miss Dina Nasr, Administrator"
$ws.Range("G107").Value = "Miss Dina Nasr, Administrator"
$ws.Range("G108").Value = "Miss Dina Nasr, Administrator"
$ws.Range("G119").Value = "Miss Dina Nasr, Administrator"
$ws.Range("G120").Value = "Miss Dina Nasr, Administrator"
$ws.Range("G122").Value = "Miss Dina Nasr, Administrator"
$ws.Range("G125").Value = "Miss Dina Nasr, Administrator"
$ws.Range("G127").Value = "Miss Dina Nasr, Administrator"
$ws.Range("G132").Value = "Miss Dina Nasr, Administrator"
$ws.Range("G133").Value = "Miss Dina Nasr, Administrator"
$ws.Range("G134").Value = "Miss Dina Nasr, Administrator"
$ws.Range("G145").Value = "Miss Dina Nasr, Administrator"
$ws.Range("G146").Value = "Miss Dina Nasr, Administrator"
$ws.Range("G148").Value = "Miss Dina Nasr, Administrator"
$ws.Range("G151").Value = "Miss Dina Nasr, Administrator"
$ws.Range("G153").Value = "Miss Dina Nasr, Administrator"
